# aircyber-v1.5.2.xlsx : v1 -> v2 Excel library layout conversion
#
# v1 layout (3 sheets):
#   library_content        -> library + framework metadata mixed together, plus "tab" rows
#   controls                -> requirement nodes (assessable/depth/ref_id/name/description/implementation_groups)
#   implementation_groups   -> ref_id/name/description triples (Bronze/Silver/Gold)
#
# v2 layout (5 sheets):
#   library_meta                   -> only library-level metadata (type/urn/version/locale/ref_id/name/description/copyright/provider/packager)
#   controls_meta                  -> framework-level metadata (type/base_urn/urn/ref_id/name/description/implementation_groups_definition)
#   controls_content                -> same requirement-node rows as old "controls" sheet, renamed
#   implementation_groups_meta      -> implementation-groups-level metadata (type/name)
#   implementation_groups_content   -> same ref_id/name/description rows as old "implementation_groups" sheet, renamed + blank cells removed
#
# NOTE on this COM host's worksheet-handle semantics: a `Worksheet` object
# captured from `Worksheets.Item(<index>)` stays bound to that *position*
# (not to the sheet's identity) for its own property/method calls once the
# collection is mutated by Add/Move/Delete. The only reliable way to keep
# addressing the right sheet after inserting new tabs is to re-resolve it
# fresh by name, via `$wb.Worksheets.Item("<name>")`, right before every use.
# So below: rename the 3 original sheets FIRST (while indices are still
# untouched), then always re-look-up by name from that point on.

$wb = $excel.ActiveWorkbook
$nl = [char]10

$description = "AirCyber is the AeroSpace and Defense official standard for Cybersecurity maturity evaluation and increase built by Airbus, Dassault Aviation, Safran and Thales to help the AeroSpace SupplyChain to be more resilient." + $nl + `
"Their joint venture BoostAeroSpace is offering this extract of the AirCyber maturity level matrix to provide further details on this standard, the questions and the AirCyber maturity levels they are associated to." + $nl + `
"AirCyber program uses this maturity level matrix as the base of the cyber maturity evaluation as is the evaluation activity is the very starting point for any cyber maturity progression. Being aware of the problems is the mandatory very first knowledge a company shall know to decide to launch a cybersecurity company program." + $nl + `
"Source: https://boostaerospace.com/aircyber/"

$copyright = [char]0xA9 + " Boost Aerospace" + $nl + `
"This work is licensed under a Creative Commons Attribution-NonCommercial-ShareAlike 4.0 International License. Any commercial use of this work must be contracted with BoostAeroSpace." + $nl + `
"Permission given to include AirCyber in CISO Assistant."

# ---------------------------------------------------------------------------
# 1) Rename the existing 3 sheets to their new v2 "content" / "meta" names
#    while the collection indices are still untouched (1,2,3).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "library_meta"
$wb.Worksheets.Item(2).Name = "controls_content"
$wb.Worksheets.Item(3).Name = "implementation_groups_content"

# ---------------------------------------------------------------------------
# 2) Insert the two brand-new "_meta" sheets at their correct tab positions
#    (always re-resolving the anchor sheet fresh, by name):
#      library_meta, controls_meta, controls_content,
#      implementation_groups_meta, implementation_groups_content
# ---------------------------------------------------------------------------
$controlsMeta = $wb.Worksheets.Add($wb.Worksheets.Item("controls_content"), $null)
$controlsMeta.Name = "controls_meta"

$igMeta = $wb.Worksheets.Add($wb.Worksheets.Item("implementation_groups_content"), $null)
$igMeta.Name = "implementation_groups_meta"

# ---------------------------------------------------------------------------
# 3) Rewrite library_meta (A1:B10) from scratch.
# ---------------------------------------------------------------------------
$libraryMeta = $wb.Worksheets.Item("library_meta")
$libraryMeta.UsedRange.ClearContents()

$libData = New-Object 'object[,]' 10,2
$libData[0,0] = "type";        $libData[0,1] = "library"
$libData[1,0] = "urn";         $libData[1,1] = "urn:intuitem:risk:library:aircyber-v1.5.2"
$libData[2,0] = "version";     $libData[2,1] = "'1"
$libData[3,0] = "locale";      $libData[3,1] = "en"
$libData[4,0] = "ref_id";      $libData[4,1] = "AirCyber-v1.5.2"
$libData[5,0] = "name";        $libData[5,1] = "Public AirCyber Maturity Level Matrix"
$libData[6,0] = "description"; $libData[6,1] = $description
$libData[7,0] = "copyright";   $libData[7,1] = $copyright
$libData[8,0] = "provider";    $libData[8,1] = "Boost Aerospace"
$libData[9,0] = "packager";    $libData[9,1] = "intuitem"
$libraryMeta.Range("A1:B10").Value = $libData

# ---------------------------------------------------------------------------
# 4) Populate the brand-new controls_meta sheet (A1:B7).
# ---------------------------------------------------------------------------
$controlsMetaData = New-Object 'object[,]' 7,2
$controlsMetaData[0,0] = "type";                              $controlsMetaData[0,1] = "framework"
$controlsMetaData[1,0] = "base_urn";                          $controlsMetaData[1,1] = "urn:intuitem:risk:req_node:aircyber-v1.5.2"
$controlsMetaData[2,0] = "urn";                               $controlsMetaData[2,1] = "urn:intuitem:risk:framework:aircyber-v1.5.2"
$controlsMetaData[3,0] = "ref_id";                            $controlsMetaData[3,1] = "AirCyber-v1.5.2"
$controlsMetaData[4,0] = "name";                              $controlsMetaData[4,1] = "Public AirCyber Maturity Level Matrix"
$controlsMetaData[5,0] = "description";                       $controlsMetaData[5,1] = $description
$controlsMetaData[6,0] = "implementation_groups_definition";  $controlsMetaData[6,1] = "implementation_groups"
$wb.Worksheets.Item("controls_meta").Range("A1:B7").Value = $controlsMetaData

# ---------------------------------------------------------------------------
# 5) Populate the brand-new implementation_groups_meta sheet (A1:B2).
# ---------------------------------------------------------------------------
$igMetaData = New-Object 'object[,]' 2,2
$igMetaData[0,0] = "type"; $igMetaData[0,1] = "implementation_groups"
$igMetaData[1,0] = "name"; $igMetaData[1,1] = "implementation_groups"
$wb.Worksheets.Item("implementation_groups_meta").Range("A1:B2").Value = $igMetaData

# ---------------------------------------------------------------------------
# 6) controls_content keeps its original requirement rows untouched (only the
#    sheet got renamed above).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 7) implementation_groups_content: drop the blank B/C cells on rows 2-4,
#    keep the ref_id/name/description header + the 3 maturity level rows.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("implementation_groups_content").Range("B2:C4").ClearContents()

Write-Output "aircyber v1 -> v2 layout conversion complete"
